$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 (Property/Value metadata table) ---
$ws1.Range("B3").Value = "6.0.0"
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"
$ws1.Range("B9").Value = "Alvearie Team"
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# Remove the old duplicate "Contact" row (row 11), shifting rows 12-21 up by one
$ws1.Rows.Item(11).Delete()

# --- Sheet2 (Extension element table) ---
$ws2.Range("K2").Value = "Match Compared To"
$ws2.Range("L2").Value = "Reference to the record that was compared for matching"
